$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue "D2" "68.781.16"
Set-TextValue "E2" "  -0.73%  "
Set-TextValue "D3" "2.455.50"
Set-TextValue "E3" "  -1.29%  "
Set-TextValue "E4" "  -0.08%  "
Set-TextValue "D5" "558.64"
Set-TextValue "E5" "  -1.60%  "
Set-TextValue "D6" "162.73"
Set-TextValue "E6" "  -1.67%  "
Set-TextValue "E7" "  -0.02%  "
Set-TextValue "D8" "0.503"
Set-TextValue "E8" "  -1.23%  "
Set-TextValue "B9" "Dogecoin"
Set-TextValue "C9" "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextValue "D9" "0.151"
Set-TextValue "E9" "  -4.33%  "
Set-TextValue "B10" "LidoStakedEther"
Set-TextValue "C10" "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
Set-TextValue "D10" "2.288.89"
Set-TextValue "E10" "  -7.98%  "
Set-TextValue "E11" "  -0.43%  "
Set-TextValue "E12" "  -3.23%  "
Set-TextValue "E13" "  -1.07%  "
Set-TextValue "D14" "2.907.43"
Set-TextValue "E14" "  -1.25%  "
Set-TextValue "D15" "68.785.35"
Set-TextValue "E15" "  -0.61%  "
Set-TextValue "E16" "  -3.08%  "
Set-TextValue "D17" "23.57"
Set-TextValue "E17" "  -2.07%  "
Set-TextValue "D18" "2.454.91"
Set-TextValue "E18" "  -1.47%  "
Set-TextValue "D19" "10.77"
Set-TextValue "E19" "  -3.36%  "
Set-TextValue "D20" "341.12"
Set-TextValue "E20" "  -3.04%  "
Set-TextValue "E21" "  -5.09%  "
Set-TextValue "D22" "3.79"
Set-TextValue "E22" "  -2.52%  "
Set-TextValue "E23" "  +1.30%  "
Set-TextValue "E24" "  +0.08%  "
Set-TextValue "D25" "66.94"
Set-TextValue "E26" "  -2.17%  "
Set-TextValue "D27" "2.582.04"
Set-TextValue "E27" "  -1.58%  "
Set-TextValue "E28" "  +0.05%  "
Set-TextValue "D29" "8.19"
Set-TextValue "E29" "  -4.87%  "
Set-TextValue "D30" "0.0₃0818"
Set-TextValue "E30" "  -5.47%  "
Set-TextValue "D31" "7.15"
Set-TextValue "E31" "  -4.41%  "
Set-TextValue "D32" "437.90"
Set-TextValue "E32" "  +0.08%  "
Set-TextValue "D33" "0.999"
Set-TextValue "E33" "  -0.08%  "
Set-TextValue "E34" "  -3.28%  "
Set-TextValue "E35" "  -5.28%  "
Set-TextValue "D36" "157.28"
Set-TextValue "E36" "  +2.62%  "
Set-TextValue "D37" "19.02"
Set-TextValue "E37" "  -0.15%  "
Set-TextValue "E38" "  +0.05%  "
Set-TextValue "E39" "  -3.55%  "
Set-TextValue "D40" "17.77"
Set-TextValue "E40" "  -1.73%  "
Set-TextValue "E41" "  -2.35%  "
Set-TextValue "E42" "  -3.72%  "
Set-TextValue "D43" "37.43"
Set-TextValue "E44" "  -5.45%  "
Set-TextValue "D45" "1.10"
Set-TextValue "E45" "  +2.96%  "
Set-TextValue "E46" "  -3.45%  "
Set-TextValue "D47" "133.02"
Set-TextValue "E47" "  -4.05%  "
Set-TextValue "E48" "  -2.09%  "
Set-TextValue "D49" "0.0718"
Set-TextValue "E49" "  -0.63%  "
Set-TextValue "E50" "  -4.00%  "
Set-TextValue "D51" "0.559"
Set-TextValue "E51" "  -2.36%  "
